$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Table "Group 272" (Key Accomplishments) ---
# Row 4, Col 2 is currently an empty cell (just an endParaRPr) -> fill it in.
$accomplTbl = $s.Shapes.Item(4).Table
$accomplCell = $accomplTbl.Cell(4, 2)
$accomplCell.Shape.TextFrame.TextRange.Text = "Learning Backend information"

# --- Table "Group 289" (Upcoming Activities) ---
# Row 1, Col 2 header text: date moves from 03/06/2019 to 03/27/2019.
$upcomingTbl = $s.Shapes.Item(5).Table
$upcomingCell = $upcomingTbl.Cell(1, 2)
$upcomingCell.Shape.TextFrame.TextRange.Text = "Upcoming Activities (Next Period through  (03/27/2019)"

# --- "Diamond 15" shape reposition ---
$diamond = $s.Shapes.Item(13)
$diamond.Left = 122.52719
$diamond.Top = 321.68925
